$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.715.44"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "1.601.04"
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("E6").Value = "  -0.71%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("D10").Value = "'19.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.75%  "
$ws.Range("D11").Value = "'0.0844"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("D12").Value = "1.825.44"
$ws.Range("E12").Value = "  +0.27%  "
$ws.Range("D13").Value = "1.611.16"
$ws.Range("E13").Value = "  +1.35%  "
$ws.Range("E14").Value = "  +0.25%  "
$ws.Range("D15").Value = "'0.523"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("D16").Value = "'64.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").Value = "26.683.45"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").Value = "0.0₃0742"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("D19").Value = "'210.17"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.80%  "
$ws.Range("E20").Value = "  +2.97%  "
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("D22").Value = "'4.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("E24").Value = "  +0.62%  "
$ws.Range("D25").Value = "'144.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.89%  "
$ws.Range("D26").Value = "'1.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("D27").Value = "'7.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  -1.00%  "
$ws.Range("D29").Value = "'15.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.64%  "
$ws.Range("E30").Value = "  +0.23%  "
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").Value = "'3.27"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.34%  "
$ws.Range("D33").Value = "'2.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.10%  "
$ws.Range("D34").Value = "1.296.29"
$ws.Range("E34").Value = "  +1.29%  "
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("E36").Value = "  +0.85%  "
$ws.Range("D37").Value = "'0.604"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.88%  "
$ws.Range("D38").Value = "'1.19"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +15.87%  "
$ws.Range("E39").Value = "  -0.60%  "
$ws.Range("E40").Value = "  -1.94%  "
$ws.Range("D41").Value = "'5.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.37%  "
$ws.Range("D42").Value = "'2.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.25%  "
$ws.Range("E43").Value = "  -0.69%  "
$ws.Range("D44").Value = "'63.11"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.67%  "
$ws.Range("D45").Value = "1.738.19"
$ws.Range("D46").Value = "'90.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.61%  "
$ws.Range("E47").Value = "  -3.01%  "
$ws.Range("E49").Value = "  +1.96%  "
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("D51").Value = "'7.38"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.88%  "
